# Update with latest cht-conf changes again and remove NO_LABEL
#
# - inputs/source and inputs/contact are now proper XLSForm "hidden" fields
#   (previously they were "text" fields with a NO_LABEL appearance, which is
#   no longer needed now that "hidden" is supported) - so the NO_LABEL
#   appearance is removed from the sheet (and therefore from sharedStrings).
# - conditionalFormatting ranges/formula for column C are tidied up: the
#   duplicate C27-only rules are merged into the main C2:C26/C28:C10000
#   ranges, and the "blank label" check also now excludes $A2="hidden" and
#   requires $H2 (calculation) to be blank.
# - the selected/active cell on the survey sheet view is reset to A2.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# --- survey sheet: row 6 ("source") becomes type=hidden, drop NO_LABEL/appearance ---
$survey.Range("A6").Value = "hidden"
$survey.Range("C6").ClearContents()
$survey.Range("F6").ClearContents()

# --- survey sheet: row 7 (begin_group contact) - drop NO_LABEL/appearance ---
$survey.Range("C7").ClearContents()

# --- survey sheet: row 9 ("name") becomes type=hidden, drop NO_LABEL/appearance ---
$survey.Range("A9").Value = "hidden"
$survey.Range("C9").ClearContents()
$survey.Range("F9").ClearContents()

# --- survey sheet: row 5 (begin_group inputs) - drop NO_LABEL/appearance ---
$survey.Range("C5").ClearContents()

# --- reset the active selection on the survey sheet back to A2 ---
$survey.Range("A2").Select()

# --- settings sheet: force the cached "now" formula to recompute ---
$settings.Range("C2").Formula = $settings.Range("C2").Formula

$excel.Calculate()
